$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value2 = 614.06665  # H33: 528.3889 -> 614.06665
$ws.Cells.Item(33, 9).Value2 = 684.25  # I33: 567.4 -> 684.25
$ws.Cells.Item(33, 11).Value2 = 684.25  # K33: 567.4 -> 684.25
$ws.Cells.Item(33, 13).Value2 = -455.25  # M33: -338.4 -> -455.25
$ws.Cells.Item(39, 8).Value2 = 385.05884  # H39: 328.65 -> 385.05884
$ws.Cells.Item(39, 9).Value2 = 196.14285  # I39: 163.11765 -> 196.14285
$ws.Cells.Item(39, 11).Value2 = 588.4285500000001  # K39: 489.35295 -> 588.4285500000001
$ws.Cells.Item(39, 13).Value2 = -292.4285500000001  # M39: -193.35295 -> -292.4285500000001
$ws.Cells.Item(48, 8).Value2 = 4551.778  # H48: 9619.333000000001 -> 4551.778
$ws.Cells.Item(48, 10).Value2 = 3421.2856  # J48: 10174.75 -> 3421.2856
$ws.Cells.Item(48, 12).Value2 = 10263.8568  # L48: 30524.25 -> 10263.8568
$ws.Cells.Item(48, 14).Value2 = -10847.8568  # N48: -31108.25 -> -10847.8568
$ws.Cells.Item(51, 8).Value2 = 18024.688  # H51: 17317.412 -> 18024.688
$ws.Cells.Item(51, 9).Value2 = 9999  # I51: 9199.200000000001 -> 9999
$ws.Cells.Item(51, 10).Value2 = 20699.916  # J51: 20700 -> 20699.916
$ws.Cells.Item(51, 11).Value2 = 9999  # K51: 9199.200000000001 -> 9999
$ws.Cells.Item(51, 12).Value2 = 20699.916  # L51: 20700 -> 20699.916
$ws.Cells.Item(51, 13).Value2 = -9515  # M51: -8715.200000000001 -> -9515
$ws.Cells.Item(51, 14).Value2 = -21667.916  # N51: -21668 -> -21667.916
$ws.Cells.Item(56, 8).Value2 = 4551.778  # H56: 9619.333000000001 -> 4551.778
$ws.Cells.Item(56, 10).Value2 = 3421.2856  # J56: 10174.75 -> 3421.2856
$ws.Cells.Item(56, 12).Value2 = 10263.8568  # L56: 30524.25 -> 10263.8568
$ws.Cells.Item(56, 14).Value2 = -11331.8568  # N56: -31592.25 -> -11331.8568
$ws.Cells.Item(75, 8).Value2 = 134081.5  # H75: 150897.8 -> 134081.5
$ws.Cells.Item(75, 10).Value2 = 59997.25  # J75: 63329.668 -> 59997.25
$ws.Cells.Item(75, 12).Value2 = 59997.25  # L75: 63329.668 -> 59997.25
$ws.Cells.Item(75, 14).Value2 = -61869.25  # N75: -65201.668 -> -61869.25
$ws.Cells.Item(78, 8).Value2 = 134081.5  # H78: 150897.8 -> 134081.5
$ws.Cells.Item(78, 10).Value2 = 59997.25  # J78: 63329.668 -> 59997.25
$ws.Cells.Item(78, 12).Value2 = 179991.75  # L78: 189989.004 -> 179991.75
$ws.Cells.Item(78, 14).Value2 = -189351.75  # N78: -199349.004 -> -189351.75
$ws.Cells.Item(107, 8).Value2 = 41669410  # H107: 38464300 -> 41669410
$ws.Cells.Item(107, 9).Value2 = 45457264  # I107: 41669410 -> 45457264
$ws.Cells.Item(107, 11).Value2 = 45457264  # K107: 41669410 -> 45457264
$ws.Cells.Item(107, 13).Value2 = -45455344  # M107: -41667490 -> -45455344
$ws.Cells.Item(116, 8).Value2 = 4100.6  # H116: 4231.6 -> 4100.6
$ws.Cells.Item(116, 9).Value2 = 3440.2  # I116: 3702.2 -> 3440.2
$ws.Cells.Item(116, 11).Value2 = 3440.2  # K116: 3702.2 -> 3440.2
$ws.Cells.Item(116, 13).Value2 = 1.800000000000182  # M116: -260.1999999999998 -> 1.800000000000182
$ws.Cells.Item(141, 8).Value2 = 0  # H141: 1750 -> 0
$ws.Cells.Item(141, 9).Value2 = 0  # I141: 1750 -> 0
$ws.Cells.Item(141, 11).Value2 = 0  # K141: 5250 -> 0
$ws.Cells.Item(141, 13).ClearContents()  # M141 was -70

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value2 = 248.33333  # H4: 266.66666 -> 248.33333
$ws.Cells.Item(4, 9).Value2 = 248.33333  # I4: 266.66666 -> 248.33333
$ws.Cells.Item(4, 11).Value2 = 248.33333  # K4: 266.66666 -> 248.33333
$ws.Cells.Item(4, 13).Value2 = -132.33333  # M4: -150.66666 -> -132.33333
$ws.Cells.Item(26, 8).Value2 = 11821  # H26: 9859.166999999999 -> 11821
$ws.Cells.Item(26, 9).Value2 = 369  # I26: 289.25 -> 369
$ws.Cells.Item(26, 11).Value2 = 369  # K26: 289.25 -> 369
$ws.Cells.Item(26, 13).Value2 = -39  # M26: 40.75 -> -39
$ws.Cells.Item(34, 8).Value2 = 204805.6  # H34: 179838 -> 204805.6
$ws.Cells.Item(34, 10).Value2 = 218342.67  # J34: 177507 -> 218342.67
$ws.Cells.Item(34, 12).Value2 = 218342.67  # L34: 177507 -> 218342.67
$ws.Cells.Item(34, 14).Value2 = -218884.67  # N34: -178049 -> -218884.67
$ws.Cells.Item(37, 8).Value2 = 45395  # H37: 24997.727 -> 45395
$ws.Cells.Item(37, 9).Value2 = 9000  # I37: 8333.333000000001 -> 9000
$ws.Cells.Item(37, 11).Value2 = 9000  # K37: 8333.333000000001 -> 9000
$ws.Cells.Item(37, 13).Value2 = -8727  # M37: -8060.333000000001 -> -8727
$ws.Cells.Item(61, 8).Value2 = 5465.364  # H61: 5180.6665 -> 5465.364
$ws.Cells.Item(61, 9).Value2 = 4305.5625  # I61: 4187.4707 -> 4305.5625
$ws.Cells.Item(61, 10).Value2 = 8558.166999999999  # J61: 7592.7144 -> 8558.166999999999
$ws.Cells.Item(61, 11).Value2 = 4305.5625  # K61: 4187.4707 -> 4305.5625
$ws.Cells.Item(61, 12).Value2 = 8558.166999999999  # L61: 7592.7144 -> 8558.166999999999
$ws.Cells.Item(61, 13).Value2 = -4093.5625  # M61: -3975.4707 -> -4093.5625
$ws.Cells.Item(61, 14).Value2 = -8982.166999999999  # N61: -8016.7144 -> -8982.166999999999
$ws.Cells.Item(74, 8).Value2 = 267226.25  # H74: 290084.22 -> 267226.25
$ws.Cells.Item(74, 9).Value2 = 668026.7  # I74: 834872.5 -> 668026.7
$ws.Cells.Item(74, 10).Value2 = 5834.696  # J74: 5846.8696 -> 5834.696
$ws.Cells.Item(74, 11).Value2 = 668026.7  # K74: 834872.5 -> 668026.7
$ws.Cells.Item(74, 12).Value2 = 5834.696  # L74: 5846.8696 -> 5834.696
$ws.Cells.Item(74, 13).Value2 = -667152.7  # M74: -833998.5 -> -667152.7
$ws.Cells.Item(74, 14).Value2 = -7582.696  # N74: -7594.8696 -> -7582.696
$ws.Cells.Item(77, 8).Value2 = 267226.25  # H77: 290084.22 -> 267226.25
$ws.Cells.Item(77, 9).Value2 = 668026.7  # I77: 834872.5 -> 668026.7
$ws.Cells.Item(77, 10).Value2 = 5834.696  # J77: 5846.8696 -> 5834.696
$ws.Cells.Item(77, 11).Value2 = 3340133.5  # K77: 4174362.5 -> 3340133.5
$ws.Cells.Item(77, 12).Value2 = 29173.48  # L77: 29234.348 -> 29173.48
$ws.Cells.Item(77, 13).Value2 = -3335765.5  # M77: -4169994.5 -> -3335765.5
$ws.Cells.Item(77, 14).Value2 = -37909.48  # N77: -37970.348 -> -37909.48
$ws.Cells.Item(132, 8).Value2 = 5578.185  # H132: 5882.9653 -> 5578.185
$ws.Cells.Item(132, 9).Value2 = 3277.2222  # I132: 3315.158 -> 3277.2222
$ws.Cells.Item(132, 10).Value2 = 10180.111  # J132: 10761.8 -> 10180.111
$ws.Cells.Item(132, 11).Value2 = 9831.6666  # K132: 9945.474 -> 9831.6666
$ws.Cells.Item(132, 12).Value2 = 30540.333  # L132: 32285.4 -> 30540.333
$ws.Cells.Item(132, 13).Value2 = -7301.6666  # M132: -7415.474 -> -7301.6666
$ws.Cells.Item(132, 14).Value2 = -35600.333  # N132: -37345.39999999999 -> -35600.333
$ws.Cells.Item(136, 8).Value2 = 5465.364  # H136: 5180.6665 -> 5465.364
$ws.Cells.Item(136, 9).Value2 = 4305.5625  # I136: 4187.4707 -> 4305.5625
$ws.Cells.Item(136, 10).Value2 = 8558.166999999999  # J136: 7592.7144 -> 8558.166999999999
$ws.Cells.Item(136, 11).Value2 = 12916.6875  # K136: 12562.4121 -> 12916.6875
$ws.Cells.Item(136, 12).Value2 = 25674.501  # L136: 22778.1432 -> 25674.501
$ws.Cells.Item(136, 13).Value2 = -10366.6875  # M136: -10012.4121 -> -10366.6875
$ws.Cells.Item(136, 14).Value2 = -30774.501  # N136: -27878.1432 -> -30774.501

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value2 = 1350.7778  # H94: 1410.92 -> 1350.7778
$ws.Cells.Item(94, 9).Value2 = 1367.8695  # I94: 1441.0952 -> 1367.8695
$ws.Cells.Item(94, 11).Value2 = 1367.8695  # K94: 1441.0952 -> 1367.8695
$ws.Cells.Item(94, 13).Value2 = -916.8695  # M94: -990.0952 -> -916.8695

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value2 = 8837.25  # H16: 7037.4287 -> 8837.25
$ws.Cells.Item(16, 9).Value2 = 0  # I16: 4750 -> 0
$ws.Cells.Item(16, 10).Value2 = 8837.25  # J16: 7952.4 -> 8837.25
$ws.Cells.Item(16, 11).Value2 = 0  # K16: 4750 -> 0
$ws.Cells.Item(16, 12).Value2 = 8837.25  # L16: 7952.4 -> 8837.25
$ws.Cells.Item(16, 13).ClearContents()  # M16 was -4463
$ws.Cells.Item(16, 14).Value2 = -9411.25  # N16: -8526.4 -> -9411.25
$ws.Cells.Item(31, 8).Value2 = 33337090  # H31: 34486570 -> 33337090
$ws.Cells.Item(31, 10).Value2 = 4916.4287  # J31: 5051.25 -> 4916.4287
$ws.Cells.Item(31, 12).Value2 = 4916.4287  # L31: 5051.25 -> 4916.4287
$ws.Cells.Item(31, 14).Value2 = -5506.4287  # N31: -5641.25 -> -5506.4287
$ws.Cells.Item(34, 8).Value2 = 33337090  # H34: 34486570 -> 33337090
$ws.Cells.Item(34, 10).Value2 = 4916.4287  # J34: 5051.25 -> 4916.4287
$ws.Cells.Item(34, 12).Value2 = 4916.4287  # L34: 5051.25 -> 4916.4287
$ws.Cells.Item(34, 14).Value2 = -5320.4287  # N34: -5455.25 -> -5320.4287
$ws.Cells.Item(113, 8).Value2 = 8837.25  # H113: 7037.4287 -> 8837.25
$ws.Cells.Item(113, 9).Value2 = 0  # I113: 4750 -> 0
$ws.Cells.Item(113, 10).Value2 = 8837.25  # J113: 7952.4 -> 8837.25
$ws.Cells.Item(113, 11).Value2 = 0  # K113: 4750 -> 0
$ws.Cells.Item(113, 12).Value2 = 8837.25  # L113: 7952.4 -> 8837.25
$ws.Cells.Item(113, 13).ClearContents()  # M113 was -2580
$ws.Cells.Item(113, 14).Value2 = -13177.25  # N113: -12292.4 -> -13177.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value2 = 53605456  # H4: 48944160 -> 53605456
$ws.Cells.Item(4, 9).Value2 = 95909430  # I4: 87917000 -> 95909430
$ws.Cells.Item(4, 10).Value2 = 7071082  # J4: 6428329 -> 7071082
$ws.Cells.Item(4, 11).Value2 = 287728290  # K4: 263751000 -> 287728290
$ws.Cells.Item(4, 12).Value2 = 21213246  # L4: 19284987 -> 21213246
$ws.Cells.Item(4, 13).Value2 = -287728178  # M4: -263750888 -> -287728178
$ws.Cells.Item(4, 14).Value2 = -21213470  # N4: -19285211 -> -21213470
$ws.Cells.Item(138, 8).Value2 = 16671167  # H138: 12504832 -> 16671167
$ws.Cells.Item(138, 9).Value2 = 50000000  # I138: 25002916 -> 50000000
$ws.Cells.Item(138, 11).Value2 = 150000000  # K138: 75008748 -> 150000000
$ws.Cells.Item(138, 13).Value2 = -149994860  # M138: -75003608 -> -149994860

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value2 = 13912.904  # H70: 13484.954 -> 13912.904
$ws.Cells.Item(70, 9).Value2 = 13895.167  # I70: 13172.308 -> 13895.167
$ws.Cells.Item(70, 11).Value2 = 13895.167  # K70: 13172.308 -> 13895.167
$ws.Cells.Item(70, 13).Value2 = -13625.167  # M70: -12902.308 -> -13625.167
$ws.Cells.Item(73, 8).Value2 = 13912.904  # H73: 13484.954 -> 13912.904
$ws.Cells.Item(73, 9).Value2 = 13895.167  # I73: 13172.308 -> 13895.167
$ws.Cells.Item(73, 11).Value2 = 13895.167  # K73: 13172.308 -> 13895.167
$ws.Cells.Item(73, 13).Value2 = -12959.167  # M73: -12236.308 -> -12959.167
$ws.Cells.Item(126, 8).Value2 = 2708.0967  # H126: 2965.8518 -> 2708.0967
$ws.Cells.Item(126, 9).Value2 = 2122.7273  # I126: 2379.2778 -> 2122.7273
$ws.Cells.Item(126, 11).Value2 = 6368.1819  # K126: 7137.8334 -> 6368.1819
$ws.Cells.Item(126, 13).Value2 = -3898.1819  # M126: -4667.8334 -> -3898.1819
$ws.Cells.Item(132, 8).Value2 = 8316.191999999999  # H132: 8445.259 -> 8316.191999999999
$ws.Cells.Item(132, 9).Value2 = 4567.4614  # I132: 4765 -> 4567.4614
$ws.Cells.Item(132, 10).Value2 = 12064.923  # J132: 11389.467 -> 12064.923
$ws.Cells.Item(132, 11).Value2 = 13702.3842  # K132: 14295 -> 13702.3842
$ws.Cells.Item(132, 12).Value2 = 36194.769  # L132: 34168.401 -> 36194.769
$ws.Cells.Item(132, 13).Value2 = -11172.3842  # M132: -11765 -> -11172.3842
$ws.Cells.Item(132, 14).Value2 = -41254.769  # N132: -39228.401 -> -41254.769

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value2 = 956.6  # H55: 1005.5714 -> 956.6
$ws.Cells.Item(55, 9).Value2 = 559.0909  # I55: 587.9 -> 559.0909
$ws.Cells.Item(55, 11).Value2 = 559.0909  # K55: 587.9 -> 559.0909
$ws.Cells.Item(55, 13).Value2 = -386.0909  # M55: -414.9 -> -386.0909
$ws.Cells.Item(61, 8).Value2 = 7666.7144  # H61: 6881.6113 -> 7666.7144
$ws.Cells.Item(61, 9).Value2 = 7640.364  # I61: 6711.6 -> 7640.364
$ws.Cells.Item(61, 10).Value2 = 7763.3335  # J61: 7731.6665 -> 7763.3335
$ws.Cells.Item(61, 11).Value2 = 7640.364  # K61: 6711.6 -> 7640.364
$ws.Cells.Item(61, 12).Value2 = 7763.3335  # L61: 7731.6665 -> 7763.3335
$ws.Cells.Item(61, 13).Value2 = -7438.364  # M61: -6509.6 -> -7438.364
$ws.Cells.Item(61, 14).Value2 = -8167.3335  # N61: -8135.6665 -> -8167.3335
$ws.Cells.Item(93, 8).Value2 = 1450  # H93: 1576.8 -> 1450
$ws.Cells.Item(93, 9).Value2 = 1450  # I93: 1500 -> 1450
$ws.Cells.Item(93, 10).Value2 = 0  # J93: 1884 -> 0
$ws.Cells.Item(93, 11).Value2 = 1450  # K93: 1500 -> 1450
$ws.Cells.Item(93, 12).Value2 = 0  # L93: 1884 -> 0
$ws.Cells.Item(93, 13).Value2 = -202  # M93: -252 -> -202
$ws.Cells.Item(93, 14).ClearContents()  # N93 was -4380
$ws.Cells.Item(113, 8).Value2 = 7666.7144  # H113: 6881.6113 -> 7666.7144
$ws.Cells.Item(113, 9).Value2 = 7640.364  # I113: 6711.6 -> 7640.364
$ws.Cells.Item(113, 10).Value2 = 7763.3335  # J113: 7731.6665 -> 7763.3335
$ws.Cells.Item(113, 11).Value2 = 7640.364  # K113: 6711.6 -> 7640.364
$ws.Cells.Item(113, 12).Value2 = 7763.3335  # L113: 7731.6665 -> 7763.3335
$ws.Cells.Item(113, 13).Value2 = -5470.364  # M113: -4541.6 -> -5470.364
$ws.Cells.Item(113, 14).Value2 = -12103.3335  # N113: -12071.6665 -> -12103.3335
$ws.Cells.Item(132, 8).Value2 = 9845  # H132: 9124.130999999999 -> 9845
$ws.Cells.Item(132, 9).Value2 = 9568.462  # I132: 9012.666999999999 -> 9568.462
$ws.Cells.Item(132, 10).Value2 = 10444.167  # J132: 9333.125 -> 10444.167
$ws.Cells.Item(132, 11).Value2 = 28705.386  # K132: 27038.001 -> 28705.386
$ws.Cells.Item(132, 12).Value2 = 31332.501  # L132: 27999.375 -> 31332.501
$ws.Cells.Item(132, 13).Value2 = -26175.386  # M132: -24508.001 -> -26175.386
$ws.Cells.Item(132, 14).Value2 = -36392.501  # N132: -33059.375 -> -36392.501

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(34, 8).Value2 = 34599.4  # H34: 34599.2 -> 34599.4
$ws.Cells.Item(34, 10).Value2 = 34749.75  # J34: 34749.5 -> 34749.75
$ws.Cells.Item(34, 12).Value2 = 34749.75  # L34: 34749.5 -> 34749.75
$ws.Cells.Item(34, 14).Value2 = -35155.75  # N34: -35155.5 -> -35155.75
$ws.Cells.Item(39, 8).Value2 = 49998  # H39: 49948 -> 49998
$ws.Cells.Item(39, 10).Value2 = 49998  # J39: 49948 -> 49998
$ws.Cells.Item(39, 12).Value2 = 49998  # L39: 49948 -> 49998
$ws.Cells.Item(39, 14).Value2 = -50824  # N39: -50774 -> -50824
$ws.Cells.Item(62, 8).Value2 = 6000  # H62: 16317.546 -> 6000
$ws.Cells.Item(62, 9).Value2 = 6000  # I62: 20187.5 -> 6000
$ws.Cells.Item(62, 10).Value2 = 0  # J62: 5997.6665 -> 0
$ws.Cells.Item(62, 11).Value2 = 6000  # K62: 20187.5 -> 6000
$ws.Cells.Item(62, 12).Value2 = 0  # L62: 5997.6665 -> 0
$ws.Cells.Item(62, 13).Value2 = -5376  # M62: -19563.5 -> -5376
$ws.Cells.Item(62, 14).ClearContents()  # N62 was -7245.6665
$ws.Cells.Item(65, 8).Value2 = 6000  # H65: 16317.546 -> 6000
$ws.Cells.Item(65, 9).Value2 = 6000  # I65: 20187.5 -> 6000
$ws.Cells.Item(65, 10).Value2 = 0  # J65: 5997.6665 -> 0
$ws.Cells.Item(65, 11).Value2 = 30000  # K65: 100937.5 -> 30000
$ws.Cells.Item(65, 12).Value2 = 0  # L65: 29988.3325 -> 0
$ws.Cells.Item(65, 13).Value2 = -26880  # M65: -97817.5 -> -26880
$ws.Cells.Item(65, 14).ClearContents()  # N65 was -36228.3325
$ws.Cells.Item(107, 8).Value2 = 2739.923  # H107: 3022.7856 -> 2739.923
$ws.Cells.Item(107, 9).Value2 = 1792.8572  # I107: 2406.25 -> 1792.8572
$ws.Cells.Item(107, 11).Value2 = 5378.571599999999  # K107: 7218.75 -> 5378.571599999999
$ws.Cells.Item(107, 13).Value2 = -3458.571599999999  # M107: -5298.75 -> -3458.571599999999
$ws.Cells.Item(113, 8).Value2 = 757.17645  # H113: 842.6875 -> 757.17645
$ws.Cells.Item(113, 9).Value2 = 596  # I113: 730.1111 -> 596
$ws.Cells.Item(113, 11).Value2 = 1788  # K113: 2190.3333 -> 1788
$ws.Cells.Item(113, 13).Value2 = 382  # M113: -20.33329999999978 -> 382
$ws.Cells.Item(122, 8).Value2 = 3974.8276  # H122: 4099.2144 -> 3974.8276
$ws.Cells.Item(122, 9).Value2 = 4033.3076  # I122: 4174.96 -> 4033.3076
$ws.Cells.Item(122, 11).Value2 = 12099.9228  # K122: 12524.88 -> 12099.9228
$ws.Cells.Item(122, 13).Value2 = -9649.9228  # M122: -10074.88 -> -9649.9228
$ws.Cells.Item(132, 8).Value2 = 4787.727  # H132: 4793.788 -> 4787.727
$ws.Cells.Item(132, 9).Value2 = 1758.7646  # I132: 1770.5294 -> 1758.7646
$ws.Cells.Item(132, 11).Value2 = 5276.293799999999  # K132: 5311.5882 -> 5276.293799999999
$ws.Cells.Item(132, 13).Value2 = -2746.293799999999  # M132: -2781.5882 -> -2746.293799999999
